$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: Insert two new columns before column D (shifts D:K -> F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Step 2: Copy number formats from column F into the two new columns D:E
# so the new cells inherit the same date / accounting number formats as
# the rest of the row (Insert() alone copies formats from the left, i.e.
# from column C, which is plain "General"). This is done per contiguous
# block of data rows (skipping the section-header rows 5, 6, 37, 79 which
# have no D:K data and must not gain any new cells).
$dataRowBlocks = @(@(7,35), @(38,77), @(80,102))
foreach ($block in $dataRowBlocks) {
    $r1 = $block[0]
    $r2 = $block[1]
    $src = $ws.Range("F$r1" + ":F$r2")
    $src.Copy()
    $dst = $ws.Range("D$r1" + ":E$r2")
    $dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}
$excel.CutCopyMode = $false

$rowDE = @{
  7 = @(43465, 43373)
  8 = @(1078000, 828000)
  9 = @(233000, 236000)
  10 = @(845000, 592000)
  11 = @($null, $null)
  12 = @(16000, 4000)
  13 = @(0, 0)
  14 = @(-31000, -2000)
  15 = @(130000, 128000)
  16 = @($null, $null)
  17 = @(582000, 641000)
  18 = @(496000, 187000)
  19 = @($null, $null)
  20 = @(-104000, -96000)
  21 = @(522000, 219000)
  22 = @(0, 0)
  23 = @(392000, 91000)
  24 = @(0, 0)
  25 = @(0, 0)
  26 = @(392000, 91000)
  27 = @(339000, 65000)
  28 = @(0, 0)
  29 = @(0, 0)
  30 = @(0, 0)
  31 = @(0, 0)
  32 = @(104000, 96000)
  33 = @(339000, 65000)
  34 = @(0, 0)
  35 = @(339000, 65000)
  38 = @(43465, 43373)
  39 = @($null, $null)
  40 = @($null, $null)
  41 = @(17000, 31000)
  42 = @(0, 0)
  43 = @(367000, 380000)
  44 = @(69000, 69000)
  45 = @(187000, 66000)
  46 = @(640000, 546000)
  47 = @(0, 0)
  48 = @(6455000, 6386000)
  49 = @(0, 0)
  50 = @(0, 0)
  51 = @(0, 0)
  52 = @(63000, 52000)
  53 = @(0, 0)
  54 = @(7158000, 6984000)
  55 = @($null, $null)
  56 = @($null, $null)
  57 = @(390000, 349000)
  58 = @("NA", "NA")
  59 = @(217000, 522000)
  60 = @(607000, 871000)
  61 = @(5251000, 5108000)
  62 = @(791000, 865000)
  63 = @(0, 0)
  64 = @(0, 0)
  65 = @(0, 0)
  66 = @(7519000, 7709000)
  67 = @($null, $null)
  68 = @(0, 0)
  69 = @(0, 0)
  70 = @(0, 0)
  71 = @(0, 0)
  72 = @(-5342000, -5688000)
  73 = @(0, 0)
  74 = @(0, 0)
  75 = @(0, 0)
  76 = @(-361000, -725000)
  77 = @(0, 0)
  80 = @(43465, 43373)
  81 = @(339000, 65000)
  82 = @($null, $null)
  83 = @(130000, 128000)
  84 = @(0, 0)
  85 = @(0, 0)
  86 = @(0, 0)
  87 = @(0, 0)
  88 = @(0, 0)
  89 = @(68000, 159000)
  90 = @($null, $null)
  91 = @(-186000, -177000)
  92 = @(0, 0)
  93 = @(0, 0)
  94 = @(-191000, -158000)
  95 = @($null, $null)
  96 = @(0, 0)
  97 = @(0, 0)
  98 = @(0, 0)
  99 = @(0, 0)
  100 = @(109000, -12000)
  101 = @(0, 0)
  102 = @(-14000, -11000)
}

# Step 3: Write the two new quarters of data (D = 2018-12-31, E = 2018-09-30)
# into every data row. Empty financial statement section rows keep $null
# (blank cell, already formatted from step 2).
foreach ($r in $rowDE.Keys) {
    $vals = $rowDE[$r]
    if ($null -ne $vals[0]) {
        $ws.Cells.Item($r, 4).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($r, 5).Value = $vals[1]
    }
}

# Step 4: Row 58 (Short Term Investments) - the two oldest periods that were
# previously shifted into F58/G58 are corrected to "NA" (matching D58/E58),
# instead of keeping the old numeric values.
$ws.Cells.Item(58, 6).Value = "NA"
$ws.Cells.Item(58, 7).Value = "NA"
